# Koopases have animations but less logic
# Adds new enemy-spawner rows (9-11, 13, 15-16) to Sheet1, matching the
# pattern already used by existing rows (e.g. row 7): columns B, C, E, F
# are computed from H/I via simple formulas, the rest are literal values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column J back down to its "auto" width (it was widened in the
# original file to fit a long value that no longer needs the extra room).
$ws.Columns(10).ColumnWidth = 8

# row number, A, G, H, I, J, K, L
$rows = @(
    @(9,  9, 301, 1328, 27,  89, 1100, 1500),
    @(10, 9, 301, 1376, 27,  89, 1100, 1500),
    @(11, 9, 301, 1425, 27,  89, 1100, 1500),
    @(13, 9, 30,  1472, 123, 89, 1160, 1480),
    @(15, 9, 31,  576,  91,  89, 515,  590),
    @(16, 9, 31,  2096, 91,  89, 2090, 2104)
)

foreach ($r in $rows) {
    $row = $r[0]

    $ws.Cells.Item($row, 1).Value = $r[1]                       # A
    $ws.Cells.Item($row, 2).Formula = "=H$row - 160"            # B
    $ws.Cells.Item($row, 3).Formula = "=I$row - 100"            # C
    $ws.Cells.Item($row, 4).Value = 0                           # D
    $ws.Cells.Item($row, 5).Formula = "= H$row + 160"           # E
    $ws.Cells.Item($row, 6).Formula = "=I$row + 16"             # F
    $ws.Cells.Item($row, 7).Value = $r[2]                       # G
    $ws.Cells.Item($row, 8).Value = $r[3]                       # H
    $ws.Cells.Item($row, 9).Value = $r[4]                       # I
    $ws.Cells.Item($row, 10).Value = $r[5]                      # J
    $ws.Cells.Item($row, 11).Value = $r[6]                      # K
    $ws.Cells.Item($row, 12).Value = $r[7]                      # L
}

# Select the newly added rows, mirroring where the author ended up after
# adding this data (matches the saved selection/active-cell state).
[void]$ws.Range("A15:L16").Select()
